# Applies the "Change prosumer output under new salvage" re-run results
# to this Grid-Search output workbook: a new grid-search pass produced
# updated dispatch/demand figures on several sheets, and the scenario's
# "Required Level of Met Demand" input flips from 0 to an unbounded ("inf")
# requirement.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B4").Value = "inf"
$ws.Range("B6").Value = 682383.9862542464
$ws.Range("B7").Value = 1783776.864811973
$ws.Range("B8").Value = 24018212.71255795
$ws.Range("B10").Value = 4558719.140590771

# ---------------------------------------------------------------
# Sheet: Costs and Revenues
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("C2").Value = 155942.7727794829
$ws.Range("D2").Value = 162628.2766561623
$ws.Range("E2").Value = 175177.03067367
$ws.Range("G2").Value = 175177.0306736701
$ws.Range("I2").Value = 175177.0306736701
$ws.Range("K2").Value = 175177.03067367
$ws.Range("P2").Value = 175177.0306736701
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 30675.49398760892
$ws.Range("E3").Value = 71144.70374179265
$ws.Range("C4").Value = 33183.28134215029
$ws.Range("D4").Value = 29187.47990269742
$ws.Range("C5").Value = 41112.97456555201
$ws.Range("D5").Value = 41893.81860702013
$ws.Range("C6").Value = 81646.51687178061
$ws.Range("D6").Value = 60871.4841588358
$ws.Range("E6").Value = 75995.72448559753

# ---------------------------------------------------------------
# Sheet: Installed Capacities
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Installed Capacities")
$ws.Range("C3").Value = 343.3658057592665
$ws.Range("D3").Value = 379.1843397715654

# ---------------------------------------------------------------
# Sheet: Added Capacities
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Added Capacities")
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 35.81853401229886
$ws.Range("E3").Value = 88.55419228182201

# ---------------------------------------------------------------
# Sheet: PV Dispatch
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G5").Value = 1.380365048278457
$ws.Range("H5").Value = 14.13666355068175
$ws.Range("I5").Value = 53.21652352375527
$ws.Range("J5").Value = 117.1567580163238
$ws.Range("K5").Value = 175.5876105099509
$ws.Range("L5").Value = 217.8319573562028
$ws.Range("M5").Value = 242.3800242835248
$ws.Range("N5").Value = 246.301986476946
$ws.Range("O5").Value = 232.5759815281271
$ws.Range("P5").Value = 198.4982193987526
$ws.Range("Q5").Value = 149.0638961072804
$ws.Range("R5").Value = 86.70935596392171
$ws.Range("S5").Value = 31.45506853764537
$ws.Range("T5").Value = 6.042547998838949
$ws.Range("U5").Value = 0.1104292038622766
$ws.Range("G6").Value = 0.7385604123878563
$ws.Range("H6").Value = 7.132938719640613
$ws.Range("I6").Value = 25.42850542651172
$ws.Range("J6").Value = 69.77776247038076
$ws.Range("K6").Value = 119.2613101003671
$ws.Range("L6").Value = 160.3615491897405
$ws.Range("M6").Value = 187.1343641388002
$ws.Range("O6").Value = 175.7223100473861
$ws.Range("P6").Value = 141.0326457655372
$ws.Range("Q6").Value = 94.276588781299
$ws.Range("R6").Value = 45.85553156913376
$ws.Range("S6").Value = 13.71843573009899
$ws.Range("T6").Value = 2.976916749931753
$ws.Range("U6").Value = 0.04858950081499056
$ws.Range("G7").Value = 0.6191842398937591
$ws.Range("H7").Value = 5.505110787419063
$ws.Range("I7").Value = 18.62055877789597
$ws.Range("J7").Value = 43.77632576048877
$ws.Range("K7").Value = 71.937950780384
$ws.Range("L7").Value = 92.05580962929581
$ws.Range("M7").Value = 97.05994407716445
$ws.Range("N7").Value = 94.75207554665141
$ws.Range("O7").Value = 87.5188778351652
$ws.Range("P7").Value = 74.88751934133245
$ws.Range("Q7").Value = 51.84823666964924
$ws.Range("R7").Value = 27.84077500467756
$ws.Range("S7").Value = 10.7906926170576
$ws.Range("T7").Value = 2.64560538863697
$ws.Range("U7").Value = 0.03377368581238691
$ws.Range("G8").Value = 1.524359154860563
$ws.Range("H8").Value = 15.61134319471575
$ws.Range("I8").Value = 58.76785631776192
$ws.Range("J8").Value = 129.3780778198468
$ws.Range("K8").Value = 193.9042008450945
$ws.Range("L8").Value = 240.5553073306585
$ws.Range("M8").Value = 267.66412945091
$ws.Range("N8").Value = 271.9952148996576
$ws.Range("O8").Value = 256.8373685535129
$ws.Range("P8").Value = 219.2047519178927
$ws.Range("Q8").Value = 164.6136396844487
$ws.Range("R8").Value = 95.75452576150994
$ws.Range("S8").Value = 34.73633424138512
$ws.Range("T8").Value = 6.672882200402118
$ws.Range("U8").Value = 0.121948732388845
$ws.Range("G9").Value = 0.8156040515841217
$ws.Range("H9").Value = 7.877018077141387
$ws.Range("I9").Value = 28.08110440761121
$ws.Range("J9").Value = 77.05669682093671
$ws.Range("K9").Value = 131.7021682770729
$ws.Range("L9").Value = 177.0898183055787
$ws.Range("M9").Value = 206.6554651755031
$ws.Range("O9").Value = 194.0529516942274
$ws.Range("P9").Value = 155.7446017450416
$ws.Range("Q9").Value = 104.1111417425626
$ws.Range("R9").Value = 50.63899541326681
$ws.Range("S9").Value = 15.14948753709979
$ws.Range("T9").Value = 3.287456681604419
$ws.Range("U9").Value = 0.05365816128842908
$ws.Range("G10").Value = 0.6837750389323309
$ws.Range("H10").Value = 6.079381709780183
$ws.Range("I10").Value = 20.56298026171047
$ws.Range("J10").Value = 48.3428952525158
$ws.Range("K10").Value = 79.4422272505017
$ws.Range("L10").Value = 101.6586998790849
$ws.Range("M10").Value = 107.1848454210016
$ws.Range("N10").Value = 104.6362293667994
$ws.Range("O10").Value = 96.64849368472623
$ws.Range("P10").Value = 82.6994828905066
$ws.Range("Q10").Value = 57.25683530550637
$ws.Range("R10").Value = 30.74501220508462
$ws.Range("S10").Value = 11.91633408757525
$ws.Range("T10").Value = 2.921584257256322
$ws.Range("U10").Value = 0.03729682030539991
$ws.Range("Q11").Value = 203.0572840332873
$ws.Range("T11").Value = 8.23125798618574
$ws.Range("U12").Value = 0.06618941491321521
$ws.Range("H13").Value = 7.499152202429724
$ws.Range("K13").Value = 97.99505638757849
$ws.Range("Q14").Value = 203.0572840332874
$ws.Range("T14").Value = 8.231257986185739
$ws.Range("U15").Value = 0.06618941491321523
$ws.Range("H16").Value = 7.499152202429723
$ws.Range("K16").Value = 97.9950563875785

# ---------------------------------------------------------------
# Sheet: Fed-in Capacity
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("M5").Value = 94.50134181136153
$ws.Range("N5").Value = 89.21461903306823
$ws.Range("O5").Value = 99.00804712831379
$ws.Range("L6").Value = 42.27423123832554
$ws.Range("M6").Value = 20.73666047215161
$ws.Range("O6").Value = 32.82469745261386
$ws.Range("P6").Value = 54.90492507792075
$ws.Range("L8").Value = 104.2530745629478
$ws.Range("M8").Value = 69.21723664397632
$ws.Range("N8").Value = 63.5213906103566
$ws.Range("O8").Value = 79.68126565070395
$ws.Range("L9").Value = 25.54596212248728
$ws.Range("M9").Value = 1.215559435448739
$ws.Range("O9").Value = 14.49405580577258
$ws.Range("P9").Value = 40.19296909841634
$ws.Range("K11").Value = 82.69308435253603
$ws.Range("K14").Value = 82.693084352536

# ---------------------------------------------------------------
# Sheet: Unmet Demand
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("G5").Value = 422.2098565259836
$ws.Range("H5").Value = 349.9597650259403
$ws.Range("I5").Value = 254.4136275388651
$ws.Range("J5").Value = 103.6577574092157
$ws.Range("K5").Value = 62.01598937483027
$ws.Range("L5").Value = 18.38450875428643
$ws.Range("O5").Value = 4.934605547775988
$ws.Range("P5").Value = 56.01237276874357
$ws.Range("Q5").Value = 120.8562777625268
$ws.Range("R5").Value = 211.4982239222711
$ws.Range("S5").Value = 235.9346815086628
$ws.Range("T5").Value = 219.4718503679457
$ws.Range("U5").Value = 248.7911763328777
$ws.Range("G6").Value = 161.0306354769367
$ws.Range("H6").Value = 137.9016867325838
$ws.Range("I6").Value = 117.4902810648562
$ws.Range("J6").Value = 82.938597341113
$ws.Range("K6").Value = 13.57083277636244
$ws.Range("Q6").Value = 55.82666100818672
$ws.Range("R6").Value = 155.3651088272201
$ws.Range("S6").Value = 208.3135435403173
$ws.Range("T6").Value = 230.1584524961191
$ws.Range("U6").Value = 249.6505853349712
$ws.Range("G7").Value = 169.247518043053
$ws.Range("H7").Value = 167.6944209426916
$ws.Range("I7").Value = 168.1305865516738
$ws.Range("J7").Value = 133.20905535072
$ws.Range("K7").Value = 88.99056004441434
$ws.Range("L7").Value = 62.96960703121103
$ws.Range("M7").Value = 59.70650310975206
$ws.Range("N7").Value = 47.31542685419205
$ws.Range("O7").Value = 72.58297639552541
$ws.Range("P7").Value = 91.43798011889524
$ws.Range("Q7").Value = 152.7386590351999
$ws.Range("R7").Value = 221.644124959155
$ws.Range("S7").Value = 243.9209403175845
$ws.Range("T7").Value = 218.1427676199599
$ws.Range("U7").Value = 291.2202965310711
$ws.Range("G8").Value = 422.0658624194015
$ws.Range("H8").Value = 348.4850853819063
$ws.Range("I8").Value = 248.8622947448584
$ws.Range("J8").Value = 91.4364376056927
$ws.Range("K8").Value = 43.69939903968668
$ws.Range("P8").Value = 35.30584024960351
$ws.Range("Q8").Value = 105.3065341853584
$ws.Range("R8").Value = 202.4530541246828
$ws.Range("S8").Value = 232.653415804923
$ws.Range("T8").Value = 218.8415161663825
$ws.Range("U8").Value = 248.7796568043511
$ws.Range("G9").Value = 160.9535918377404
$ws.Range("H9").Value = 137.157607375083
$ws.Range("I9").Value = 114.8376820837567
$ws.Range("J9").Value = 75.65966299055705
$ws.Range("K9").Value = 1.129974599656691
$ws.Range("Q9").Value = 45.99210804692308
$ws.Range("R9").Value = 150.5816449830871
$ws.Range("S9").Value = 206.8824917333164
$ws.Range("T9").Value = 229.8479125644464
$ws.Range("U9").Value = 249.6455166744977
$ws.Range("G10").Value = 169.1829272440145
$ws.Range("H10").Value = 167.1201500203304
$ws.Range("I10").Value = 166.1881650678593
$ws.Range("J10").Value = 128.642485858693
$ws.Range("K10").Value = 81.48628357429664
$ws.Range("L10").Value = 53.36671678142191
$ws.Range("M10").Value = 49.58160176591487
$ws.Range("N10").Value = 37.43127303404407
$ws.Range("O10").Value = 63.45336054596437
$ws.Range("P10").Value = 83.62601656972109
$ws.Range("Q10").Value = 147.3300603993428
$ws.Range("R10").Value = 218.7398877587479
$ws.Range("S10").Value = 242.7952988470669
$ws.Range("T10").Value = 217.8667887513406
$ws.Range("U10").Value = 291.2167733965781
$ws.Range("Q11").Value = 66.86288983651983
$ws.Range("K13").Value = 62.93345443721985
$ws.Range("L13").Value = 29.62548300147246
$ws.Range("Q14").Value = 66.8628898365198
$ws.Range("K16").Value = 62.93345443721984
$ws.Range("L16").Value = 29.62548300147245

# ---------------------------------------------------------------
# Sheet: Household Surplus
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B3").Value = 373422.0583155675
$ws.Range("B4").Value = 384012.8107888695

